$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 values
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "ssss"
$ws.Range("C3").Value = "ssss"
$ws.Range("D3").Value = "ssss"
$ws.Range("E3").Value = "ssss"
$ws.Range("F3").Value = "ssss"
$ws.Range("G3").Value = "ssss"

# A3 should carry the same formatting as A2 (bold/bordered/centered style)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
